# Dividend Calculation workbook update
# - December (row 14) Taxable Account dividend on the "Yearly" sheet increased
#   from 57.19 to 133.62. All the SUM()/cross-sheet formulas that depend on
#   it (Yearly!G14, Yearly!D15, Yearly!G15, 'All Time'!F7, 'All Time'!I7,
#   'All Time'!F46, 'All Time'!I46) recalculate automatically.
# - Restore the selection/active-cell state that Excel persisted after the edit.

$wb = $excel.ActiveWorkbook

$wsYearly  = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# --- the actual data edit -------------------------------------------------
$wsYearly.Range("D14").Value = 133.62

# --- recalculate everything so cached formula results are fresh ----------
$excel.Calculate()

# --- restore on-screen selection state ------------------------------------
$wsYearly.Activate() | Out-Null
$wsYearly.Range("M9").Select() | Out-Null

$wsAllTime.Activate() | Out-Null
$wsAllTime.Range("P26").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
